# Weekly refresh of the Ají (Hortaliza) daily-price rows: each data row
# (2-16, 18; row 17 is untouched) is re-stamped with another row's
# date/variety/quality/volume/price/unit data, per the updated source feed.
# Columns A,B,C,E,F,G,O,R are identical for every row, so only
# D,H,I,J,K,L,M,N,P,Q actually move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> source row it should copy its D/H/I/J/K/L/M/N/P/Q values from.
$rowMap = @{
    2  = 13
    3  = 14
    4  = 15
    5  = 4
    6  = 9
    7  = 10
    8  = 12
    9  = 11
    10 = 18
    11 = 16
    12 = 5
    13 = 8
    14 = 3
    15 = 7
    16 = 6
    18 = 2
}

$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot the "before" values for every row/column we might read from,
# so writes to earlier rows don't corrupt reads for later rows.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    if (-not $snapshot.ContainsKey($src)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Range("$c$src").Value2
        }
        $snapshot[$src] = $rowVals
    }
}

foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}
